# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting so values
# like "1.001" or "6.500" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.243.34"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.906.95"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "307.19"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.5249"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.3807"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").Value = "0.07289"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "21.76"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").Value = "0.9028"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "0.08183"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "96.16"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.459.21"
$ws.Range("E15").Value = "  -23.52%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "0.000008674"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "27.279.24"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "5.114"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "10.84"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Value = "6.500"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "150.12"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").Value = "2.322"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "1.743"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "116.73"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").Value = "4.848"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "4.843"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "0.09235"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "0.8348"
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("D33").Value = "0.05052"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "1.229"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "2.995"
$ws.Range("D36").Value = "2.740"
$ws.Range("E36").Value = "  +5.50%  "
$ws.Range("D37").Value = "3.340"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").Value = "0.5808"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "0.02009"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "9.170"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "6.609"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "0.4943"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "10.23"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D48").Value = "1.644"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "38.89"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "64.48"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "0.06062"
$ws.Range("E51").Value = "  +1.72%  "
